$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '65.713.08'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.95%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.450.78'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -3.78%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '595.94'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.59%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '136.61'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -7.40%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.452.86'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -3.73%  '
$ws.Range("E8").Value = '  -0.12%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.489'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.05%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.54'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -4.24%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.122'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -9.99%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.380'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -8.01%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.032.29'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -3.85%  '
$ws.Range("E14").Value = '  -11.50%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '26.51'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -10.12%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.454.58'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -3.88%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '65.601.33'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.17%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.115'
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '9.91'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -10.65%  '
$ws.Range("E20").Value = '  -8.72%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.78'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -6.92%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '393.39'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -6.74%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.547'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -10.31%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '73.18'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -6.43%  '
$ws.Range("E25").Value = '  -0.04%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.591.42'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -3.71%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0000106'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -11.63%  '
$ws.Range("E28").Value = '  -0.69%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.33'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -10.82%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.26'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -9.06%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.19'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -12.47%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.455.03'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -3.59%  '
$ws.Range("E33").Value = '  +0.01%  '
$ws.Range("E34").Value = '  -6.45%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '22.84'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -8.70%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '172.93'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.95%  '
$ws.Range("E37").Value = '  -13.22%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.88'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -10.94%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.52'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -8.40%  '
$ws.Range("E40").Value = '  -13.35%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0777'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -8.62%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.816'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -7.14%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '43.58'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -5.30%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.999'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.00%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '4.44'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -14.11%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.63'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -12.29%  '
$ws.Range("B47").Value = 'EnergySwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '23.13'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.02%  '
$ws.Range("B48").Value = 'ONDO'
$ws.Range("C48").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.12'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.57%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '6.53'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -8.31%  '
$ws.Range("E50").Value = '  -15.86%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.207.55'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -7.47%  '
